{"js": "// The original document has three trailing paragraphs whose paragraph\n// properties (<w:pPr>) carry an explicit \"justify\" alignment\n// (<w:jc w:val=\"both\"/>). The edit removes that direct formatting so the\n// paragraphs fall back to the default (left) alignment, leaving their\n// <w:pPr> empty.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/alignment\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.alignment === Word.Alignment.justified) {\n    paragraph.alignment = Word.Alignment.left;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The original document has three trailing paragraphs whose paragraph\n# properties (<w:pPr>) carry an explicit \"justify\" alignment\n# (<w:jc w:val=\"both\"/>). The edit removes that direct formatting so the\n# paragraphs fall back to the default (left) alignment, leaving their\n# <w:pPr> empty.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.ParagraphFormat.Alignment -eq 3) {   # wdAlignParagraphJustify\n        $p.Range.ParagraphFormat.Alignment = 0        # wdAlignParagraphLeft\n    }\n}\n"}
